$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet "data": add new date column AI with % values (row 1-58)
# and bump the footnote date in row 59.
# ============================================================
$ws1 = $wb.Worksheets.Item("data")

# Header cell AI1: copy style (bold / centered / bordered) from AH1,
# then overwrite with the new date label.
$ws1.Range("AH1").Copy($ws1.Range("AI1"))
$ws1.Range("AI1").Value = "28. 9. 2021"

# New data column AI (rows 2-58)
$ws1.Range("AI2").Value = 0.07000000000000001
$ws1.Range("AI3").Value = 0.03
$ws1.Range("AI4").Value = 0.16
$ws1.Range("AI5").Value = 0.1
$ws1.Range("AI6").Value = 0.06
$ws1.Range("AI7").Value = 0.18
$ws1.Range("AI8").Value = 0.02
$ws1.Range("AI9").Value = 0.02
$ws1.Range("AI10").Value = 0.21
$ws1.Range("AI11").Value = 0.09
$ws1.Range("AI12").Value = 0.015
$ws1.Range("AI13").Value = 0.08
$ws1.Range("AI14").Value = 0.09
$ws1.Range("AI15").Value = 0.03
$ws1.Range("AI16").Value = 0.17
$ws1.Range("AI17").Value = 0.04
$ws1.Range("AI18").Value = 0.02
$ws1.Range("AI19").Value = 0.14
$ws1.Range("AI20").Value = 0.06
$ws1.Range("AI21").Value = 0.03
$ws1.Range("AI22").Value = 0.18
$ws1.Range("AI23").Value = 0.09
$ws1.Range("AI24").Value = 0.06
$ws1.Range("AI25").Value = 0.22
$ws1.Range("AI26").Value = 0.06
$ws1.Range("AI27").Value = 0.01
$ws1.Range("AI28").Value = 0.08
$ws1.Range("AI29").Value = 0.07000000000000001
$ws1.Range("AI30").Value = 0.05
$ws1.Range("AI31").Value = 0.18
$ws1.Range("AI32").Value = 0.06
$ws1.Range("AI33").Value = 0.04
$ws1.Range("AI34").Value = 0.19
$ws1.Range("AI35").Value = 0.09
$ws1.Range("AI36").Value = 0
$ws1.Range("AI37").Value = 0.34
$ws1.Range("AI38").Value = 0.04
$ws1.Range("AI39").Value = 0
$ws1.Range("AI40").Value = 0.16
$ws1.Range("AI41").Value = 0.16
$ws1.Range("AI42").Value = 0.04
$ws1.Range("AI43").Value = 0.02
$ws1.Range("AI44").Value = 0.04
$ws1.Range("AI45").Value = 0.04
$ws1.Range("AI46").Value = 0.18
$ws1.Range("AI47").Value = 0.13
$ws1.Range("AI48").Value = 0
$ws1.Range("AI49").Value = 0.19
$ws1.Range("AI50").Value = 0.05
$ws1.Range("AI51").Value = 0.05
$ws1.Range("AI52").Value = 0.24
$ws1.Range("AI53").Value = 0.11
$ws1.Range("AI54").Value = 0.02
$ws1.Range("AI55").Value = 0.11
$ws1.Range("AI56").Value = 0.05
$ws1.Range("AI57").Value = 0.07000000000000001
$ws1.Range("AI58").Value = 0.24

# Update footnote text in row 59 (new "aktualizace" date)
$ws1.Range("A59").Value = "Život během pandemie, Testování, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# ============================================================
# Sheet "pocetR": add new date column AH with sample-size values
# (row 1-20) and bump the footnote date in row 21.
# ============================================================
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AH1: copy style from AG1, then overwrite with the new date label.
$ws2.Range("AG1").Copy($ws2.Range("AH1"))
$ws2.Range("AH1").Value = "28. 9. 2021"

# New data column AH (rows 2-20)
$ws2.Range("AH2").Value = 1855
$ws2.Range("AH3").Value = 456
$ws2.Range("AH4").Value = 678
$ws2.Range("AH5").Value = 721
$ws2.Range("AH6").Value = 897
$ws2.Range("AH7").Value = 958
$ws2.Range("AH8").Value = 883
$ws2.Range("AH9").Value = 160
$ws2.Range("AH10").Value = 568
$ws2.Range("AH11").Value = 244
$ws2.Range("AH12").Value = 833
$ws2.Range("AH13").Value = 65
$ws2.Range("AH14").Value = 72
$ws2.Range("AH15").Value = 73
$ws2.Range("AH16").Value = 337
$ws2.Range("AH17").Value = 110
$ws2.Range("AH18").Value = 288
$ws2.Range("AH19").Value = 157
$ws2.Range("AH20").Value = 103

# Row 21: blank placeholder cell under the new AH column (matches the
# existing pattern of empty cells B21:AG21 on that row), then update the footnote.
$ws2.Range("AG21").Copy($ws2.Range("AH21"))
$ws2.Range("A21").Value = "Život během pandemie, Testování, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"
